$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (EmptySheet), so it
# lands at the end of the tab strip as sheet index 9 / sheetId 9.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "WithOptionalFields"

# Row 4: headers (written left-to-right so new shared strings land in the
# same order as the source commit).
$newSheet.Range("D4").Value = "Name"
$newSheet.Range("E4").Value = "Middle Name"
$newSheet.Range("F4").Value = "Surname"
$newSheet.Range("G4").Value = "UpdatedDate"

# Row 5: first optional-fields record.
$newSheet.Range("D5").Value = "Aydin"
$newSheet.Range("E5").Value = "N/A"
$newSheet.Range("F5").Value = "Eraydin"

# Row 6: second record - only some of the optional columns populated.
$newSheet.Range("D6").Value = "ABC"
$newSheet.Range("F6").Value = "DEF"
$newSheet.Range("G6").Value = "2019-02-01"

# Age column added last (matches shared-string ordering from the diff).
$newSheet.Range("H4").Value = "Age"
$newSheet.Range("H5").Value = 30

# Autosize the UpdatedDate column like Excel did when the sheet was authored.
$newSheet.Columns.Item(7).AutoFit() | Out-Null

# Make this new sheet the active tab/selection, matching the saved view.
$newSheet.Range("H6").Select() | Out-Null

Write-Output "Added WithOptionalFields sheet"
